# Update Name of Algo
# Applies updated imputed values (RandomForest result data) to specific cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.073799999999991
$ws.Range("A9").Value = -20.43539999999997
$ws.Range("B11").Value = 5.567999999999996
$ws.Range("A18").Value = -23.05410000000001
$ws.Range("A20").Value = -22.22830000000002
$ws.Range("D21").Value = -7.699699999999996
